$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 21:52"

# --- Re-rank country rows whose shared-string label swapped position ---
# Set from the bottom row up so the reused strings do not collide while Excel
# reassigns/dedupes the shared-string table.

# Sudafrica overtakes Catar (row 50 <-> row 51)
$ws.Range("A51").Value = "Catar"
$ws.Range("A50").Value = "Sudafrica"

# Costa Rica overtakes Republica de Chipre (row 85 <-> row 86)
$ws.Range("A86").Value = "Republica de Chipre"
$ws.Range("A85").Value = "Costa Rica"

# Uzbekistan overtakes Cuba (row 94 <-> row 95)
$ws.Range("A95").Value = "Cuba"
$ws.Range("A94").Value = "Uzbekistan"

# Tanzania jumps ahead of Benin/Birmania/Gabon/Haiti (rows 156-160 rotate)
$ws.Range("A160").Value = "Haiti"
$ws.Range("A159").Value = "Gabon"
$ws.Range("A158").Value = "Birmania"
$ws.Range("A157").Value = "Benin"
$ws.Range("A156").Value = "Tanzania"

# --- Updated case/death statistics ---
$ws.Range("B4").Value = 331519
$ws.Range("C4").Value = 20162
$ws.Range("D4").Value = 17115
$ws.Range("E4").Value = 304920
$ws.Range("F4").Value = 8573
$ws.Range("G4").Value = 1032
$ws.Range("H4").Value = 9484

$ws.Range("F15").Value = 1385

$ws.Range("B16").Value = 15425
$ws.Range("C16").Value = 1513
$ws.Range("D16").Value = 2847
$ws.Range("E16").Value = 12301

$ws.Range("B19").Value = 10627
$ws.Range("C19").Value = 267
$ws.Range("E19").Value = 10040
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 460

$ws.Range("B39").Value = 2402
$ws.Range("C39").Value = 223
$ws.Range("E39").Value = 1880

$ws.Range("B50").Value = 1655
$ws.Range("C50").Value = 70
$ws.Range("D50").Value = 95
$ws.Range("E50").Value = 1549
$ws.Range("F50").Value = 7
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 11

$ws.Range("B51").Value = 1604
$ws.Range("C51").Value = 279
$ws.Range("D51").Value = 123
$ws.Range("E51").Value = 1477
$ws.Range("F51").Value = 37
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 4

$ws.Range("B57").Value = 1308
$ws.Range("C57").Value = 83
$ws.Range("D57").Value = 28
$ws.Range("E57").Value = 1243
$ws.Range("G57").Value = 5
$ws.Range("H57").Value = 37

$ws.Range("B71").Value = 700
$ws.Range("C71").Value = 12
$ws.Range("D71").Value = 431
$ws.Range("E71").Value = 265

$ws.Range("B85").Value = 454
$ws.Range("C85").Value = 19
$ws.Range("D85").Value = 16
$ws.Range("E85").Value = 436
$ws.Range("F85").Value = 14
$ws.Range("H85").Value = 2

$ws.Range("B86").Value = 446
$ws.Range("C86").Value = 20
$ws.Range("D86").Value = 37
$ws.Range("E86").Value = 400
$ws.Range("F86").Value = 11
$ws.Range("H86").Value = 9

$ws.Range("B94").Value = 342
$ws.Range("C94").Value = 76
$ws.Range("D94").Value = 30
$ws.Range("E94").Value = 310
$ws.Range("F94").Value = 8
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 2

$ws.Range("B95").Value = 320
$ws.Range("C95").Value = 32
$ws.Range("D95").Value = 15
$ws.Range("E95").Value = 297
$ws.Range("F95").Value = 11
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 8

$ws.Range("B139").Value = 52
$ws.Range("C139").Value = 4
$ws.Range("E139").Value = 52

$ws.Range("C156").Value = 2
$ws.Range("D156").Value = 3
$ws.Range("E156").Value = 18
$ws.Range("H156").Value = 1

$ws.Range("B157").Value = 22
$ws.Range("C157").Value = 6
$ws.Range("D157").Value = 5
$ws.Range("E157").Value = 17
$ws.Range("H157").Value = 0

$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 20

$ws.Range("C159").Value = 0
$ws.Range("G159").Value = 0

$ws.Range("B160").Value = 21
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 1
$ws.Range("E160").Value = 19
$ws.Range("G160").Value = 1
